$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected fixed-recourse data (objective, gap, solve time) for instances 1-10
$data = @{
    2  = @{ B = -105.20922885539122;  C = 0.09054680418354322;   D = 129.891672026 }
    3  = @{ B = -100.25578112906734;  C = 0.06162279836591759;   D = 138.498554965 }
    4  = @{ B = -103.67140169613268;  C = 0.02259346017021363;   D = 226.030132661 }
    5  = @{ B = -103.09145033800306;  C = 0.09806946180741519;   D = 72.251829062  }
    6  = @{ B = -102.05315995415073;  C = 0.09107443008400061;   D = 73.510411201  }
    7  = @{ B = -102.01057739411644;  C = 0.0957130934643146;    D = 130.960530718 }
    8  = @{ B = -97.44343274822958;   C = 0.006862235569463699;  D = 11.926138511  }
    9  = @{ B = -102.53029754612697;  C = 0.09900204050256638;   D = 134.412641328 }
    10 = @{ B = -102.29121992425158;  C = 0.06533247987843727;   D = 50.896578923  }
    11 = @{ B = -99.01915304554905;   C = 0.08516894509857384;   D = 118.647207056 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
}

$wb.Save()
